$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.901.27"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "2.419.11"
$ws.Range("E3").Value = "  -1.43%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  -0.96%  "
$ws.Range("D9").Value = "2.403.53"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.340"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "60.861.98"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "2.408.59"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.88"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "578.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").Value = "2.532.61"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").Value = "0.0₃0920"
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.15%  "
$ws.Range("E33").Value = "  -3.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.39%  "
$ws.Range("E37").Value = "  -3.38%  "
$ws.Range("E38").Value = "  -2.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.82%  "
$ws.Range("E44").Value = "  -3.58%  "
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").Value = "0.0₆0283"
$ws.Range("E46").Value = "  +14.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0505"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.37%  "
